$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: A1 = "relation", B1 = "count"
$ws.Range("A1").Value = "relation"
$ws.Range("B1").Value = "count"

# Match the selection state recorded in the saved file
$ws.Range("I10").Select()
